# Revert "Merging 0.1.8 w VitalSigns"
#
# This undoes the VitalSigns merge that had been applied to the
# Observations sheet:
#   - the CIBMTR Observation Laboratory Results profile row (row 32) goes
#     back to pointing at the generic FHIR observation-codes ValueSet
#     instead of the us-core laboratory-test-codes ValueSet
#   - the "Observation Category Codes#laboratory" category label reverts
#     to its pre-merge "null#laboratory" placeholder text
#   - the three Vital Signs rows (cibmtr-vital-signs,
#     cibmtr-vital-signs-height, cibmtr-vital-signs-weight) that were
#     appended as rows 33-35 are removed again

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the three trailing Vital Signs rows (33-35) entirely.
$ws.Range("A33:K35").EntireRow.Delete()

# Restore the category-code label used by the cytogenetics and
# observation-lab profile rows (C2, C31, C32 all share this text).
$ws.Range("C2").Value2 = "null#laboratory"
$ws.Range("C31").Value2 = "null#laboratory"
$ws.Range("C32").Value2 = "null#laboratory"

# Restore the CIBMTR Observation Laboratory Results Profile's Code VS
# back to the generic FHIR observation-codes ValueSet.
$ws.Range("F32").Value2 = "http://hl7.org/fhir/ValueSet/observation-codes (extensible)"
